$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 6
$ws.Cells.Item(2, 2).Value = "Cannot track subject No image Please retry ."
$ws.Cells.Item(2, 3).Value = "No image"
$ws.Cells.Item(2, 4).Value = "3-4"
$ws.Cells.Item(2, 5).Value = "Missing"

$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "Cannot track subject No image Please retry ."
$ws.Cells.Item(3, 3).Value = "Please retry"
$ws.Cells.Item(3, 4).Value = "5-6"
$ws.Cells.Item(3, 5).Value = "Missing"

$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Cannot track subject No image Please retry ."
$ws.Cells.Item(4, 3).Value = "No image Please retry"
$ws.Cells.Item(4, 4).Value = "3-6"
$ws.Cells.Item(4, 5).Value = "'False"
$ws.Cells.Item(4, 5).Style = "Normal"

$ws.Cells.Item(5, 1).Value = 23
$ws.Cells.Item(5, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(5, 3).Value = "Compass Interference"
$ws.Cells.Item(5, 4).Value = "0-1"
$ws.Cells.Item(5, 5).Value = "Missing"

$ws.Cells.Item(6, 1).Value = 23
$ws.Cells.Item(6, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(6, 3).Value = "Temp Max Altitude: nnn"
$ws.Cells.Item(6, 4).Value = "2-5"
$ws.Cells.Item(6, 5).Value = "Missing"

$ws.Cells.Item(7, 1).Value = 23
$ws.Cells.Item(7, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(7, 3).Value = "Compass Interference Temp"
$ws.Cells.Item(7, 4).Value = "0-2"
$ws.Cells.Item(7, 5).Value = "'False"
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(8, 1).Value = 23
$ws.Cells.Item(8, 2).Value = "Compass Interference Temp Max Altitude: nnn ."
$ws.Cells.Item(8, 3).Value = "Max Altitude: nnn"
$ws.Cells.Item(8, 4).Value = "3-5"
$ws.Cells.Item(8, 5).Value = "'False"
$ws.Cells.Item(8, 5).Style = "Normal"

$ws.Cells.Item(9, 1).Value = 24
$ws.Cells.Item(9, 2).Value = "No SD card Insert card ."
$ws.Cells.Item(9, 3).Value = "No SD card"
$ws.Cells.Item(9, 4).Value = "0-2"
$ws.Cells.Item(9, 5).Value = "Missing"

$ws.Cells.Item(10, 1).Value = 24
$ws.Cells.Item(10, 2).Value = "No SD card Insert card ."
$ws.Cells.Item(10, 3).Value = "Insert card"
$ws.Cells.Item(10, 4).Value = "3-4"
$ws.Cells.Item(10, 5).Value = "Missing"

$ws.Cells.Item(11, 1).Value = 24
$ws.Cells.Item(11, 2).Value = "No SD card Insert card ."
$ws.Cells.Item(11, 3).Value = "No SD card Insert card"
$ws.Cells.Item(11, 4).Value = "0-4"
$ws.Cells.Item(11, 5).Value = "'False"
$ws.Cells.Item(11, 5).Style = "Normal"

$ws.Cells.Item(12, 1).Value = 31
$ws.Cells.Item(12, 2).Value = "Warning: Battery Temperature Below 15°C (59F) Warm battery to above 25°C (77F) before flying ."
$ws.Cells.Item(12, 3).Value = "Warning: Battery Temperature Below 15°C (59F)"
$ws.Cells.Item(12, 4).Value = "0-5"
$ws.Cells.Item(12, 5).Value = "Missing"

$ws.Cells.Item(13, 1).Value = 31
$ws.Cells.Item(13, 2).Value = "Warning: Battery Temperature Below 15°C (59F) Warm battery to above 25°C (77F) before flying ."
$ws.Cells.Item(13, 3).Value = "Warm battery to above 25°C (77F) before flying"
$ws.Cells.Item(13, 4).Value = "6-13"
$ws.Cells.Item(13, 5).Value = "Missing"

$ws.Cells.Item(14, 1).Value = 31
$ws.Cells.Item(14, 2).Value = "Warning: Battery Temperature Below 15°C (59F) Warm battery to above 25°C (77F) before flying ."
$ws.Cells.Item(14, 3).Value = "Below 15°C (59F)"
$ws.Cells.Item(14, 4).Value = "3-5"
$ws.Cells.Item(14, 5).Value = "'False"
$ws.Cells.Item(14, 5).Style = "Normal"

$ws.Cells.Item(15, 1).Value = 31
$ws.Cells.Item(15, 2).Value = "Warning: Battery Temperature Below 15°C (59F) Warm battery to above 25°C (77F) before flying ."
$ws.Cells.Item(15, 3).Value = "Warm battery to above 25°C (77F)"
$ws.Cells.Item(15, 4).Value = "6-11"
$ws.Cells.Item(15, 5).Value = "'False"
$ws.Cells.Item(15, 5).Style = "Normal"

$ws.Cells.Item(16, 1).Value = 32
$ws.Cells.Item(16, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(16, 3).Value = "Compass Interference"
$ws.Cells.Item(16, 4).Value = "0-1"
$ws.Cells.Item(16, 5).Value = "Missing"

$ws.Cells.Item(17, 1).Value = 32
$ws.Cells.Item(17, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(17, 3).Value = "Temp Max Altitude 98ft"
$ws.Cells.Item(17, 4).Value = "2-5"
$ws.Cells.Item(17, 5).Value = "Missing"

$ws.Cells.Item(18, 1).Value = 32
$ws.Cells.Item(18, 2).Value = "Compass Interference Temp Max Altitude 98ft ."
$ws.Cells.Item(18, 3).Value = "Compass Interference Temp"
$ws.Cells.Item(18, 4).Value = "0-2"
$ws.Cells.Item(18, 5).Value = "'False"
$ws.Cells.Item(18, 5).Style = "Normal"

$ws.Cells.Item(19, 1).Value = 36
$ws.Cells.Item(19, 2).Value = "Incompatible firmware version Go to Profile > Settings to update firmware ."
$ws.Cells.Item(19, 3).Value = "Go to Profile > Settings to update firmware"
$ws.Cells.Item(19, 4).Value = "3-10"
$ws.Cells.Item(19, 5).Value = "Missing"

$ws.Cells.Item(20, 1).Value = 43
$ws.Cells.Item(20, 2).Value = "Cannot change current flight mode Enable Multiple Flight Modes in Settings menu to change ."
$ws.Cells.Item(20, 3).Value = "Enable Multiple Flight Modes in Settings menu to change"
$ws.Cells.Item(20, 4).Value = "5-13"
$ws.Cells.Item(20, 5).Value = "Missing"

$ws.Cells.Item(21, 1).Value = 43
$ws.Cells.Item(21, 2).Value = "Cannot change current flight mode Enable Multiple Flight Modes in Settings menu to change ."
$ws.Cells.Item(21, 3).Value = "Multiple Flight Modes in Settings menu to change"
$ws.Cells.Item(21, 4).Value = "6-13"
$ws.Cells.Item(21, 5).Value = "'False"
$ws.Cells.Item(21, 5).Style = "Normal"
